$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.85
$ws.Range("G3").Value = 2.02
$ws.Range("P3").Value = 2.36
$ws.Range("AF3").Value = 15
$ws.Range("AK3").Value = 25
$ws.Range("AL3").Value = 48
$ws.Range("AN3").Value = 11.5
$ws.Range("Q4").Value = 1.51
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 6.6
$ws.Range("F6").Value = 2.46
$ws.Range("G6").Value = 2.5
$ws.Range("I6").Value = 3.35
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 65
$ws.Range("AE6").Value = 95
$ws.Range("AI6").Value = 55
$ws.Range("AN6").Value = 24
$ws.Range("AO6").Value = 100
$ws.Range("F7").Value = 2.6
$ws.Range("G7").Value = 2.66
$ws.Range("H7").Value = 2.62
$ws.Range("I7").Value = 2.66
$ws.Range("N7").Value = 6.4
$ws.Range("T7").Value = 1.49
$ws.Range("U7").Value = 2.9
$ws.Range("AA7").Value = 40
$ws.Range("AE7").Value = 24
$ws.Range("AF7").Value = 23
$ws.Range("AG7").Value = 12.5
$ws.Range("AI7").Value = 28
$ws.Range("F8").Value = 9.4
$ws.Range("G8").Value = 9.8
$ws.Range("I8").Value = 1.4
$ws.Range("J8").Value = 5.5
$ws.Range("K8").Value = 5.7
$ws.Range("S8").Value = 2.58
$ws.Range("U8").Value = 1.98
$ws.Range("Z8").Value = 9
$ws.Range("AE8").Value = 14
$ws.Range("AF8").Value = 95
$ws.Range("H9").Value = 19
$ws.Range("I9").Value = 20
$ws.Range("J9").Value = 8.6
$ws.Range("K9").Value = 8.8
$ws.Range("R9").Value = 1.76
$ws.Range("X9").Value = 40
$ws.Range("AC9").Value = 19.5
$ws.Range("AH9").Value = 1000
$ws.Range("F10").Value = 1.3
$ws.Range("K10").Value = 7
$ws.Range("Q10").Value = 1.39
$ws.Range("R10").Value = 2.02
$ws.Range("T10").Value = 1.73
$ws.Range("U10").Value = 2.28
$ws.Range("X10").Value = 48
$ws.Range("AC10").Value = 17
$ws.Range("AD10").Value = 42
$ws.Range("AF10").Value = 11
$ws.Range("AK10").Value = 12.5
$ws.Range("AL10").Value = 28
$ws.Range("AM10").Value = 95
$ws.Range("AN10").Value = 3.5
$ws.Range("G11").Value = 6.2
$ws.Range("H11").Value = 1.62
$ws.Range("I11").Value = 1.63
$ws.Range("J11").Value = 4.4
$ws.Range("N11").Value = 4.4
$ws.Range("AE11").Value = 16.5
$ws.Range("AG11").Value = 24
$ws.Range("AH11").Value = 22
$ws.Range("AJ11").Value = 210
$ws.Range("AK11").Value = 90
$ws.Range("I12").Value = 2.46
$ws.Range("P12").Value = 2.24
$ws.Range("Q12").Value = 1.76
$ws.Range("U12").Value = 2.46
$ws.Range("F13").Value = 2.24
$ws.Range("G13").Value = 2.5
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 4.1
$ws.Range("J13").Value = 3.25
$ws.Range("K13").Value = 3.75
$ws.Range("P13").Value = 1.83
$ws.Range("Q13").Value = 1.97
